$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 11.41231320823865
$ws.Range("E2").Value = 10.95168045806885
$ws.Range("F2").Value = 12.40086274295785
$ws.Range("G2").Value = 10.76535665477439
$ws.Range("H2").Value = 1299391231
$ws.Range("I2").Value = "IFX GR"
$ws.Range("AL2").Value = "EUR"

# Row 3
$ws.Range("D3").Value = 11.62969021187249
$ws.Range("E3").Value = 10.56350612831116
$ws.Range("F3").Value = 11.96093262605668
$ws.Range("G3").Value = 10.18671831406672
$ws.Range("H3").Value = 1299391231
$ws.Range("I3").Value = "IFX GR"
$ws.Range("AL3").Value = "EUR"

# Row 4
$ws.Range("D4").Value = 10.56350560646126
$ws.Range("E4").Value = 11.58828446102142
$ws.Range("F4").Value = 12.38015956891461
$ws.Range("G4").Value = 10.06353783079802
$ws.Range("H4").Value = 1299391231
$ws.Range("I4").Value = "IFX GR"
$ws.Range("AL4").Value = "EUR"

# Row 5
$ws.Range("D5").Value = 13.97943758394593
$ws.Range("E5").Value = 12.71140232086181
$ws.Range("F5").Value = 13.97943758394593
$ws.Range("G5").Value = 11.69179912915183
$ws.Range("H5").Value = 1299391231
$ws.Range("I5").Value = "IFX GR"
$ws.Range("AL5").Value = "EUR"

# Row 6
$ws.Range("I6").Value = "IFX GR"
$ws.Range("AL6").Value = "EUR"

# Row 7
$ws.Range("D7").Value = 13.7956901867409
$ws.Range("E7").Value = 15.59192159080505
$ws.Range("F7").Value = 15.91850857079651
$ws.Range("G7").Value = 12.25756817160668
$ws.Range("H7").Value = 1299391231
$ws.Range("I7").Value = "IFX GR"
$ws.Range("AL7").Value = "EUR"

# Row 8
$ws.Range("D8").Value = 16.72971009620683
$ws.Range("E8").Value = 17.23012588119507
$ws.Range("F8").Value = 17.35127908132555
$ws.Range("G8").Value = 16.32937686539401
$ws.Range("H8").Value = 1299391231
$ws.Range("I8").Value = "IFX GR"
$ws.Range("AL8").Value = "EUR"

# Row 9
$ws.Range("D9").Value = 17.33021137288823
$ws.Range("E9").Value = 17.8991047782898
$ws.Range("F9").Value = 18.16775050508547
$ws.Range("G9").Value = 16.9246104588649
$ws.Range("H9").Value = 1299391231
$ws.Range("I9").Value = "IFX GR"
$ws.Range("AL9").Value = "EUR"

# Row 10
$ws.Range("D10").Value = 20.52060018753459
$ws.Range("E10").Value = 20.2751644744873
$ws.Range("F10").Value = 20.55261433785486
$ws.Range("G10").Value = 19.37878861910236
$ws.Range("H10").Value = 1299391231
$ws.Range("I10").Value = "IFX GR"
$ws.Range("AL10").Value = "EUR"

# Row 11
$ws.Range("D11").Value = 19.78962598112524
$ws.Range("E11").Value = 19.61888793182373
$ws.Range("F11").Value = 21.00613509123801
$ws.Range("G11").Value = 19.02663943283996
$ws.Range("H11").Value = 1299391231
$ws.Range("I11").Value = "IFX GR"
$ws.Range("AL11").Value = "EUR"

# Row 12
$ws.Range("D12").Value = 22.87892256675036
$ws.Range("E12").Value = 25.082511302948
$ws.Range("F12").Value = 25.19989448316438
$ws.Range("G12").Value = 22.86291447422016
$ws.Range("H12").Value = 1299391231
$ws.Range("I12").Value = "IFX GR"
$ws.Range("AL12").Value = "EUR"

# Row 13
$ws.Range("I13").Value = "IFX GR"
$ws.Range("AL13").Value = "EUR"

# Row 14
$ws.Range("D14").Value = 23.46208176838126
$ws.Range("E14").Value = 22.9656449584961
$ws.Range("F14").Value = 24.05564993575423
$ws.Range("G14").Value = 21.86484760818376
$ws.Range("H14").Value = 1299391231
$ws.Range("I14").Value = "IFX GR"
$ws.Range("AL14").Value = "EUR"

# Row 15
$ws.Range("D15").Value = 22.98722863627405
$ws.Range("E15").Value = 24.45495849990845
$ws.Range("F15").Value = 25.1564463111685
$ws.Range("G15").Value = 22.39366250176973
$ws.Range("H15").Value = 1299391231
$ws.Range("I15").Value = "IFX GR"
$ws.Range("AL15").Value = "EUR"

# Row 16
$ws.Range("D16").Value = 21.14717098583352
$ws.Range("E16").Value = 19.1128549079895
$ws.Range("F16").Value = 22.4368318225283
$ws.Range("G16").Value = 17.00299530374785
$ws.Range("H16").Value = 1299391231
$ws.Range("I16").Value = "IFX GR"
$ws.Range("AL16").Value = "EUR"

# Row 17
$ws.Range("D17").Value = 18.74052815651449
$ws.Range("E17").Value = 20.95291541099548
$ws.Range("F17").Value = 21.56266910012841
$ws.Range("G17").Value = 17.64512793242761
$ws.Range("H17").Value = 1299391231
$ws.Range("I17").Value = "IFX GR"
$ws.Range("AL17").Value = "EUR"

# Row 18
$ws.Range("D18").Value = 19.69639573316351
$ws.Range("E18").Value = 23.00648404312133
$ws.Range("F18").Value = 23.65208829245456
$ws.Range("G18").Value = 19.69201907357391
$ws.Range("H18").Value = 1299391231
$ws.Range("I18").Value = "IFX GR"
$ws.Range("AL18").Value = "EUR"

# Row 19
$ws.Range("D19").Value = 18.16445449294981
$ws.Range("E19").Value = 18.62403714942932
$ws.Range("F19").Value = 19.95026102890931
$ws.Range("G19").Value = 17.10741396562588
$ws.Range("H19").Value = 1299391231
$ws.Range("I19").Value = "IFX GR"
$ws.Range("AL19").Value = "EUR"

# Row 20
$ws.Range("D20").Value = 18.23886101922508
$ws.Range("E20").Value = 19.0157754650116
$ws.Range("F20").Value = 20.11877295359612
$ws.Range("G20").Value = 16.71786229120623
$ws.Range("H20").Value = 1299391231
$ws.Range("I20").Value = "IFX GR"
$ws.Range("AL20").Value = "EUR"

# Row 21
$ws.Range("D21").Value = 22.22409967481839
$ws.Range("E21").Value = 21.42092523002625
$ws.Range("F21").Value = 24.18279784373782
$ws.Range("G21").Value = 21.25459963994453
$ws.Range("H21").Value = 1299391231
$ws.Range("I21").Value = "IFX GR"
$ws.Range("AL21").Value = "EUR"

# Row 22
$ws.Range("I22").Value = "IFX GR"
$ws.Range("AL22").Value = "EUR"

# Row 23
$ws.Range("D23").Value = 22.85280379887807
$ws.Range("E23").Value = 23.4123503074646
$ws.Range("F23").Value = 26.1491480422388
$ws.Range("G23").Value = 22.52040022533345
$ws.Range("H23").Value = 1299391231
$ws.Range("I23").Value = "IFX GR"
$ws.Range("AL23").Value = "EUR"

# Row 24
$ws.Range("D24").Value = 27.47876628737924
$ws.Range("E24").Value = 26.48155338668823
$ws.Range("F24").Value = 31.39005930539134
$ws.Range("G24").Value = 25.89984603739936
$ws.Range("H24").Value = 1299391231
$ws.Range("I24").Value = "IFX GR"
$ws.Range("AL24").Value = "EUR"

# Row 25
$ws.Range("D25").Value = 34.78058337482976
$ws.Range("E25").Value = 36.70298970031738
$ws.Range("F25").Value = 39.79988899110095
$ws.Range("G25").Value = 34.38169948989559
$ws.Range("H25").Value = 1299391231
$ws.Range("I25").Value = "IFX GR"
$ws.Range("AL25").Value = "EUR"

# Row 26
$ws.Range("D26").Value = 40.79614969212006
$ws.Range("E26").Value = 37.42341357421875
$ws.Range("F26").Value = 41.5933416334489
$ws.Range("G26").Value = 36.69312037561483
$ws.Range("H26").Value = 1299391231
$ws.Range("I26").Value = "IFX GR"
$ws.Range("AL26").Value = "EUR"

# Row 27
$ws.Range("D27").Value = 38.01991829505475
$ws.Range("E27").Value = 35.82903417205811
$ws.Range("F27").Value = 38.16486346247938
$ws.Range("G27").Value = 33.68832517568221
$ws.Range("H27").Value = 1299391231
$ws.Range("I27").Value = "IFX GR"
$ws.Range("AL27").Value = "EUR"

# Row 28
$ws.Range("D28").Value = 39.00665525031088
$ws.Range("E28").Value = 45.03855450057983
$ws.Range("F28").Value = 45.08315366863616
$ws.Range("G28").Value = 37.86382602482929
$ws.Range("H28").Value = 1299391231
$ws.Range("I28").Value = "IFX GR"
$ws.Range("AL28").Value = "EUR"

# Row 29
$ws.Range("D29").Value = 45.54585767104963
$ws.Range("E29").Value = 40.51184231567383
$ws.Range("F29").Value = 46.17580924366443
$ws.Range("G29").Value = 37.93630194959927
$ws.Range("H29").Value = 1299391231
$ws.Range("I29").Value = "IFX GR"
$ws.Range("AL29").Value = "EUR"

# Row 30
$ws.Range("D30").Value = 34.871513162342
$ws.Range("E30").Value = 30.86323010635376
$ws.Range("F30").Value = 35.13011379862758
$ws.Range("G30").Value = 29.36223041026563
$ws.Range("H30").Value = 1299391231
$ws.Range("I30").Value = "IFX GR"
$ws.Range("AL30").Value = "EUR"

# Row 31
$ws.Range("D31").Value = 25.12907906347629
$ws.Range("E31").Value = 29.85132074165345
$ws.Range("F31").Value = 30.37414067004134
$ws.Range("G31").Value = 23.24580230862804
$ws.Range("H31").Value = 1299391231
$ws.Range("I31").Value = "IFX GR"
$ws.Range("AL31").Value = "EUR"

# Row 32
$ws.Range("D32").Value = 25.11221179238521
$ws.Range("E32").Value = 27.70382398033142
$ws.Range("F32").Value = 29.70515624206894
$ws.Range("G32").Value = 24.8255051349176
$ws.Range("H32").Value = 1299391231
$ws.Range("I32").Value = "IFX GR"
$ws.Range("AL32").Value = "EUR"

# Row 33
$ws.Range("D33").Value = 32.19557972357918
$ws.Range("E33").Value = 37.00777117538452
$ws.Range("F33").Value = 38.34573853529102
$ws.Range("G33").Value = 32.16184864305057
$ws.Range("H33").Value = 1299391231
$ws.Range("I33").Value = "IFX GR"
$ws.Range("AL33").Value = "EUR"

# Row 34
$ws.Range("D34").Value = 42.42663731619185
$ws.Range("E34").Value = 37.32749711608887
$ws.Range("F34").Value = 42.70456459966367
$ws.Range("G34").Value = 35.74500794121225
$ws.Range("H34").Value = 1299391231
$ws.Range("I34").Value = "IFX GR"
$ws.Range("AL34").Value = "EUR"

# Row 35
$ws.Range("D35").Value = 42.99383383152298
$ws.Range("E35").Value = 45.37607609558106
$ws.Range("F35").Value = 45.68236512851421
$ws.Range("G35").Value = 39.26732375526938
$ws.Range("H35").Value = 1299391231
$ws.Range("I35").Value = "IFX GR"
$ws.Range("AL35").Value = "EUR"

# Row 36
$ws.Range("D36").Value = 35.91516532524682
$ws.Range("E36").Value = 31.16769380187988
$ws.Range("F36").Value = 37.74722636637873
$ws.Range("G36").Value = 30.70826023625083
$ws.Range("H36").Value = 1299391231
$ws.Range("I36").Value = "IFX GR"
$ws.Range("AL36").Value = "EUR"

# Row 37
$ws.Range("D37").Value = 42.88038990294108
$ws.Range("E37").Value = 38.28038996505737
$ws.Range("F37").Value = 43.13563292062151
$ws.Range("G37").Value = 36.71491849717615
$ws.Range("H37").Value = 1299391231
$ws.Range("I37").Value = "IFX GR"
$ws.Range("AL37").Value = "EUR"

# Row 38
$ws.Range("D38").Value = 36.1671431326616
$ws.Range("E38").Value = 37.51431176376343
$ws.Range("F38").Value = 39.48633243963496
$ws.Range("G38").Value = 34.00594253294112
$ws.Range("H38").Value = 1299391231
$ws.Range("I38").Value = "IFX GR"
$ws.Range("AL38").Value = "EUR"

# Row 39
$ws.Range("D39").Value = 39.97933828197714
$ws.Range("E39").Value = 36.74613710021973
$ws.Range("F39").Value = 41.57300784070154
$ws.Range("G39").Value = 34.5218794650166
$ws.Range("H39").Value = 1299391231
$ws.Range("I39").Value = "IFX GR"
$ws.Range("AL39").Value = "EUR"

# Row 40
$ws.Range("I40").Value = "IFX GR"
$ws.Range("AL40").Value = "EUR"

# Row 41
$ws.Range("D41").Value = 36.00089548272408
$ws.Range("E41").Value = 36.75187037277222
$ws.Range("F41").Value = 39.56658782725973
$ws.Range("G41").Value = 35.00341965478231
$ws.Range("H41").Value = 1299391231
$ws.Range("I41").Value = "IFX GR"
$ws.Range("AL41").Value = "EUR"

# Row 42
$ws.Range("D42").Value = 35.83228920555114
$ws.Range("E42").Value = 33.42572929382325
$ws.Range("F42").Value = 36.13310947036743
$ws.Range("G42").Value = 26.8076900882721
$ws.Range("H42").Value = 1299391231
$ws.Range("I42").Value = "IFX GR"
$ws.Range("AL42").Value = "EUR"

# Row 43
$ws.Range("D43").Value = 41.46109394073486
$ws.Range("E43").Value = 40.07269305801392
$ws.Range("F43").Value = 44.77590088272095
$ws.Range("G43").Value = 39.08345858764648
$ws.Range("H43").Value = 1299391231
$ws.Range("I43").Value = "IFX GR"
$ws.Range("AL43").Value = "EUR"
